$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices + 1h volume deltas); rows 40-44, 47-48, 51 also
# reshuffle coin identity (name/link) as the ranking changed.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.459.23'
$ws.Range("E2").Value = '  -2.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.198.01'
$ws.Range("E3").Value = '  -3.47%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.38'
$ws.Range("E5").Value = '  -1.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.01'
$ws.Range("E6").Value = '  -5.24%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.199.28'
$ws.Range("E8").Value = '  -3.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.506'
$ws.Range("E9").Value = '  -3.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").Value = '  -4.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.30'
$ws.Range("E11").Value = '  -3.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  -3.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000237'
$ws.Range("E13").Value = '  -5.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.52'
$ws.Range("E14").Value = '  -4.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.725.24'
$ws.Range("E15").Value = '  -3.41%  '

$ws.Range("E16").Value = '  -0.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.203.48'
$ws.Range("E17").Value = '  -3.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.562.53'
$ws.Range("E18").Value = '  -2.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.65'
$ws.Range("E19").Value = '  -3.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.91'
$ws.Range("E20").Value = '  -4.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.03'
$ws.Range("E21").Value = '  -1.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.711'
$ws.Range("E22").Value = '  -4.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.68'
$ws.Range("E23").Value = '  -4.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.42'
$ws.Range("E24").Value = '  -1.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.59'
$ws.Range("E25").Value = '  -1.18%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.71'
$ws.Range("E27").Value = '  -2.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.91'
$ws.Range("E29").Value = '  -3.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("E30").Value = '  -5.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.07'
$ws.Range("E31").Value = '  -4.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.36'
$ws.Range("E32").Value = '  -5.00%  '

$ws.Range("E33").Value = '  -3.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.43'
$ws.Range("E34").Value = '  -5.31%  '

$ws.Range("E35").Value = '  -5.79%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.86'
$ws.Range("E36").Value = '  -2.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.34'
$ws.Range("E37").Value = '  -3.69%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0696'
$ws.Range("E38").Value = '  -8.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0388'
$ws.Range("E39").Value = '  -3.37%  '

$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.987.71'
$ws.Range("E40").Value = '  -1.44%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '413.42'
$ws.Range("E41").Value = '  -4.80%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("E42").Value = '  -4.19%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.114'
$ws.Range("E43").Value = '  +2.96%  '

$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.09'
$ws.Range("E44").Value = '  -4.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.17'
$ws.Range("E45").Value = '  -2.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.252'
$ws.Range("E46").Value = '  -6.55%  '

$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.84'
$ws.Range("E48").Value = '  +1.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.89'
$ws.Range("E49").Value = '  -2.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.30'
$ws.Range("E50").Value = '  +0.93%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.112'
$ws.Range("E51").Value = '  -3.17%  '
